$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1307.95
$ws.Range("D2").Value = 4638.95
$ws.Range("E2").Value = 1765.61

$ws.Range("C3").Value = 206.19
$ws.Range("D3").Value = 697.52
$ws.Range("E3").Value = 971.0855

$ws.Range("C4").Value = 42.31
$ws.Range("D4").Value = 162.59
$ws.Range("E4").Value = 176.561

$ws.Range("C5").Value = 40.18
$ws.Range("D5").Value = 144.5
$ws.Range("E5").Value = 264.8415

$ws.Range("C6").Value = 23.37
$ws.Range("D6").Value = 52.49

$ws.Range("C7").Value = 46.93
$ws.Range("D7").Value = 773.48

$ws.Range("C8").Value = 199.15
$ws.Range("D8").Value = 536.14

$ws.Range("C9").Value = 0.75
$ws.Range("D9").Value = 1.17
$ws.Range("E9").Value = 17.6561

$ws.Range("C10").Value = 11.46
$ws.Range("D10").Value = 49.11
$ws.Range("E10").Value = 176.561

$ws.Range("C11").Value = 7.2
$ws.Range("D11").Value = 32.98
$ws.Range("E11").Value = 105.9366

$ws.Range("C12").Value = 1411.22
$ws.Range("D12").Value = 2808.04

$ws.Range("C13").Value = 1418.91
$ws.Range("D13").Value = 5221.63

$ws.Range("C14").Value = 7.3
$ws.Range("D14").Value = 33.44

$ws.Range("C15").Value = 175.35
$ws.Range("D15").Value = 710.01

$ws.Range("C16").Value = 0.44
$ws.Range("D16").Value = 5.96

$ws.Range("C17").Value = 0.99
$ws.Range("D17").Value = 6.22

$ws.Range("C18").Value = 0.54
$ws.Range("D18").Value = 6.04

$ws.Range("C19").Value = 5.9
$ws.Range("D19").Value = 41.81

$ws.Range("C20").Value = 1.54
$ws.Range("D20").Value = 32.38

$ws.Range("C21").Value = 6.48
$ws.Range("D21").Value = 86.90000000000001

$ws.Range("C22").Value = 481.54
$ws.Range("D22").Value = 10106.89

$ws.Range("C23").Value = 1.32
$ws.Range("D23").Value = 9.609999999999999

$ws.Range("C24").Value = 192.1
$ws.Range("D24").Value = 435.65

$ws.Range("C25").Value = 511.21
$ws.Range("D25").Value = 2425.36

$ws.Range("C26").Value = 6.86
$ws.Range("D26").Value = 22.04
